$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (lambda = none): refreshed statistic columns B-K
$ws.Range("B2").Value = 0.790655727213881
$ws.Range("C2").Value = 0.96479924671636
$ws.Range("D2").Value = 0.269411938726706
$ws.Range("E2").Value = 0.825586938163414
$ws.Range("F2").Value = 0.300712399516157
$ws.Range("G2").Value = 0.0039356663201051
$ws.Range("H2").Value = 0.0017076901041137
$ws.Range("I2").Value = 0.00292702740346671
$ws.Range("J2").Value = 0.000983084654728542
$ws.Range("K2").Value = 0.00217451300354179

# Row 3: new lambda value (must stay text, like the rest of column A) plus
# refreshed statistic columns B-K
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "0.0153912087055105"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = 0.790094439344204
$ws.Range("C3").Value = 0.974901470212843
$ws.Range("D3").Value = 0.227495187295058
$ws.Range("E3").Value = 0.825275184115689
$ws.Range("F3").Value = 0.270553779878389
$ws.Range("G3").Value = 0.00400898541089515
$ws.Range("H3").Value = 0.000983865440256281
$ws.Range("I3").Value = 0.00660686629302962
$ws.Range("J3").Value = 0.000751217251980834
$ws.Range("K3").Value = 0.0063557987711503

# Row 4: lambda label unchanged, refreshed statistic columns B-K
$ws.Range("B4").Value = 0.790069987873999
$ws.Range("C4").Value = 0.974262626261916
$ws.Range("D4").Value = 0.229527975063751
$ws.Range("E4").Value = 0.825171249150638
$ws.Range("F4").Value = 0.271810481025635
$ws.Range("G4").Value = 0.00405568814001805
$ws.Range("H4").Value = 0.00208663499661409
$ws.Range("I4").Value = 0.00502521396002505
$ws.Range("J4").Value = 0.00126978429789727
$ws.Range("K4").Value = 0.00458481059562871

# Row 5: lambda label unchanged, refreshed statistic columns
$ws.Range("B5").Value = 0.76769247023468
$ws.Range("D5").Value = 0.000302805589923109
$ws.Range("E5").Value = 0.799866636000343
$ws.Range("F5").Value = 0.000484235993190166
$ws.Range("G5").Value = 0.00439481325137891
$ws.Range("I5").Value = 0.000246622865204638
$ws.Range("J5").Value = 0.0000499728925881218
$ws.Range("K5").Value = 0.000394346542243373
